$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.747.81'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.455.94'
$ws.Range("E3").Value = '  -2.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.20'
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.32'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -3.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.453.68'
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("E10").Value = '  -3.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.155'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  -3.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.15'
$ws.Range("E14").Value = '  -2.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.898.67'
$ws.Range("E15").Value = '  -2.47%  '
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.395.97'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.455.92'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.94'
$ws.Range("E19").Value = '  -4.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'
$ws.Range("E20").Value = '  -3.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.18'
$ws.Range("E21").Value = '  -3.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '319.98'
$ws.Range("E22").Value = '  -2.48%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.90'
$ws.Range("E24").Value = '  +4.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.79'
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0986'
$ws.Range("E26").Value = '  -6.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.574.98'
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '535.64'
$ws.Range("E29").Value = '  -3.09%  '
$ws.Range("E30").Value = '  -5.09%  '
$ws.Range("E31").Value = '  -5.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.73'
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("E33").Value = '  -5.18%  '
$ws.Range("E34").Value = '  -3.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.74'
$ws.Range("E36").Value = '  -5.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.80'
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.26'
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.66'
$ws.Range("E42").Value = '  -7.18%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.25'
$ws.Range("E44").Value = '  -1.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.28'
$ws.Range("E45").Value = '  -3.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '143.14'
$ws.Range("E46").Value = '  -5.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.97'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("E48").Value = '  -2.71%  '
$ws.Range("E49").Value = '  -4.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.590'
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0933'
$ws.Range("E51").Value = '  -2.99%  '
